$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "한글 자연어 처리 패키지가 포함된 딥러닝 전용 도커 배포 (손쉬운 설치 및 실행)"
$ws.Range("E4").Value = "https://teddylee777.github.io/linux/docker-kaggle-ko2"

$ws.Range("D46").Value = "[LG화학] 2022년 07월, 생물정보학(Bioinformatics 채용), 기반기술연구소 R&D 경력사원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/480"

$ws.Range("D51").Value = "[python] selenium 크롤링 find_element_by_css_selector 더 이상 사용 불가"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-selenium-%ED%81%AC%EB%A1%A4%EB%A7%81-findelementbycssselector-%EB%8D%94-%EC%9D%B4%EC%83%81-%EC%82%AC%EC%9A%A9-%EB%B6%88%EA%B0%80"
